# Recruitment.xlsx edit: add a "Login" test-data sheet (positioned between
# EditCandidate and SearchVacancy), tweak EditCandidate's column widths /
# active-tab state, and make AddVacancy the active tab with an updated
# D2 value + selection.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Login" worksheet right before "SearchVacancy" (i.e.
#    right after "EditCandidate") so the final sheet order is:
#    AddCandidate, EditVacancy, EditCandidate, Login, SearchVacancy,
#    AddVacancy, SearchCandidate.
# ---------------------------------------------------------------------
$searchVacancy = $wb.Worksheets.Item("SearchVacancy")
$login = $wb.Worksheets.Add($searchVacancy)
$login.Name = "Login"

# ---------------------------------------------------------------------
# 2. Populate the Login sheet's test data. Cells are written in the
#    specific order below so that new shared strings are interned in the
#    same sequence as the target workbook.
# ---------------------------------------------------------------------
$login.Range("A1").Value = 'username'
$login.Range("B1").Value = 'password'
$login.Range("A3").Value = 'Admin'
$login.Range("B3").Value = 'admin123'
$login.Range("A5").Value = 'ghh'
$login.Range("B7").Value = 'ddd'
$login.Range("D3").Value = 'login success'
$login.Range("C3").Value = 'valid'
$login.Range("C2").Value = 'invalid'
$login.Range("D1").Value = 'alertmsg'
$login.Range("D4").Value = 'Invalid credentials'
$login.Range("D2").Value = 'Username cannot be empty'
$login.Range("D5").Value = 'Password cannot be empty'
$login.Range("C1").Value = 'status'
$login.Range("A4").Value = 'admin123'
$login.Range("B4").Value = 'Admin'
$login.Range("C4").Value = 'invalid'
$login.Range("C5").Value = 'invalid'
$login.Range("B6").Value = 878
$login.Range("C6").Value = 'invalid'
$login.Range("D6").Value = 'Username cannot be empty'
$login.Range("A7").Value = 'Admin'
$login.Range("C7").Value = 'invalid'
$login.Range("D7").Value = 'Invalid credentials'
$login.Range("A8").Value = 'Admin'
$login.Range("B8").Value = 'admin123'
$login.Range("C8").Value = 'valid'
$login.Range("D8").Value = 'login success'

# Column widths (character units) chosen so the persisted <col> widths land
# on the target values (11.7109375 / 10.42578125 / 17.28515625 / 29.140625).
$login.Columns.Item(1).ColumnWidth = 10.833333333333334
$login.Columns.Item(2).ColumnWidth = 9.666666666666666
$login.Columns.Item(3).ColumnWidth = 16.5
$login.Columns.Item(4).ColumnWidth = 28.333333333333332

# D2/D4/D5/D6/D7 pick up the same "alert message" style already used
# elsewhere in the workbook (AddVacancy!B2 -- Consolas 9pt grey font).
$styleSource = $wb.Worksheets.Item("AddVacancy").Range("B2")
$styleSource.Copy()
foreach ($addr in @("D2", "D4", "D5", "D6", "D7")) {
    $login.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$login.PageSetup.PaperSize = 9
$login.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 3. EditCandidate: add a width for the (new) 9th column.
# ---------------------------------------------------------------------
$editCandidate = $wb.Worksheets.Item("EditCandidate")
$editCandidate.Columns.Item(9).ColumnWidth = 18.333333333333332

# ---------------------------------------------------------------------
# 4. Login sheet keeps its own selection (D9) but is not the active tab.
# ---------------------------------------------------------------------
$login.Activate()
$login.Range("D9").Select()

# ---------------------------------------------------------------------
# 5. AddVacancy becomes the active/selected tab, D2 changes 3 -> 6, and
#    the selection moves to B2.
# ---------------------------------------------------------------------
$addVacancy = $wb.Worksheets.Item("AddVacancy")
$addVacancy.Range("D2").Value = 6
$addVacancy.Activate()
$addVacancy.Range("B2").Select()

# Tab-bar / tab-scroll ratio tweak (best effort -- not all hosts persist this).
$excel.ActiveWindow.TabRatio = 0.846
